# Append 45 new master-data rows (102-146) to the
# "master-reg_center_machine_devic" sheet, continuing the existing
# regcntr_id / machine_id / device_id pattern, and leave the selection
# on the first empty row below the data (mirrors what the author did
# in Excel before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id / machine_id cycle repeated every 9 rows
$idCycle = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$startRow = 102
$endRow = 146
$deviceId = 3000121

for ($row = $startRow; $row -le $endRow; $row++) {
    $pair = $idCycle[($row - $startRow) % 9]

    $ws.Cells.Item($row, 1).Value = $pair[0]      # A: regcntr_id
    $ws.Cells.Item($row, 2).Value = $pair[1]      # B: machine_id
    $ws.Cells.Item($row, 3).Value = $deviceId     # C: device_id
    $ws.Cells.Item($row, 4).Value = "eng"         # D: lang_code
    $ws.Cells.Item($row, 5).Value = $true         # E: is_active
    $ws.Cells.Item($row, 6).Value = "superadmin()" # F: cr_by
    $ws.Cells.Item($row, 7).Value = "now()"       # G: cr_dtimes
    $ws.Cells.Item($row, 8).Value = "now()"       # H: eff_dtimes

    $deviceId++
}

# Select the row just below the new data, matching the state the
# workbook was left in when it was last saved.
$ws.Range("A147:XFD1048576").Select()
